# Update the date line and all math expressions in the table.
$d = $word.ActiveDocument

$replacements = @(
    @("2025-05-25 Sunday", "2025-05-26 Monday"),
    @("333×6=", "479×5="),
    @("576×5=", "913×9="),
    @("568×4=", "590×8="),
    @("647×8=", "781×3="),
    @("173×3=", "812×7="),
    @("362×6=", "215×3="),
    @("537×4=", "887×3="),
    @("524×5=", "809×2="),
    @("704×4=", "597×3="),
    @("967×3=", "559×6="),
    @("300×2=", "806×8="),
    @("302×5=", "154×9="),
    @("992×4=", "358×6="),
    @("234×2=", "638×6="),
    @("877×9=", "299×3="),
    @("960×6=", "711×9="),
    @("403×7=", "934×9="),
    @("252×8=", "564×6="),
    @("172×9=", "203×2="),
    @("703×2=", "574×4="),
    @("321×9=", "819×2="),
    @("156×6=", "685×7="),
    @("967×2=", "225×4="),
    @("783×5=", "839×7="),
    @("104×2=", "906×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
